# Update cryptocurrency price/volume data for cryptos.xlsx
# Commit: Updated cryptos list on Sat Sep 21 23:16:24 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.184.91"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "'2.575.47"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'588.13"
$ws.Range("E5").Value = "  +3.54%  "
$ws.Range("D6").Value = "'148.50"
$ws.Range("E6").Value = "  +0.99%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").Value = "'0.110"
$ws.Range("E9").Value = "  +3.88%  "
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +0.08%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "'27.65"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "'3.039.81"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").Value = "'63.162.02"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "'0.0000149"
$ws.Range("E16").Value = "  +3.49%  "
$ws.Range("D17").Value = "'2.576.39"
$ws.Range("E17").Value = "  +1.04%  "
$ws.Range("D18").Value = "'11.39"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'4.45"
$ws.Range("E19").Value = "  +2.93%  "
$ws.Range("D20").Value = "'343.42"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -3.60%  "
$ws.Range("D24").Value = "'66.73"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "'2.669.97"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.171"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'8.25"
$ws.Range("E28").Value = "  +11.98%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "'8.50"
$ws.Range("E29").Value = "  +0.55%  "
$ws.Range("B30").Value = "SuiNetwork"
$ws.Range("C30").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D30").Value = "'1.49"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "'1.99"
$ws.Range("E32").Value = "  +7.23%  "
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").Value = "'465.41"
$ws.Range("E34").Value = "  +12.71%  "
$ws.Range("D35").Value = "'1.64"
$ws.Range("E35").Value = "  +4.27%  "
$ws.Range("D36").Value = "'176.58"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").Value = "'0.408"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").Value = "'19.25"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").Value = "'4.65"
$ws.Range("E39").Value = "  +5.70%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'151.64"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("D44").Value = "'3.85"
$ws.Range("E44").Value = "  +1.90%  "
$ws.Range("D45").Value = "'21.11"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("E46").Value = "  +5.50%  "
$ws.Range("D47").Value = "'0.614"
$ws.Range("E47").Value = "  +1.56%  "
$ws.Range("D48").Value = "'0.0977"
$ws.Range("E48").Value = "  +1.42%  "
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("E51").Value = "  +0.67%  "
